$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.477.96"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "1.913.70"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.25"
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4835"
$ws.Range("E7").Value = "  +2.55%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4072"
$ws.Range("E8").Value = "  +1.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08168"
$ws.Range("E9").Value = "  +2.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.015"
$ws.Range("E10").Value = "  +2.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.70"
$ws.Range("E11").Value = "  +5.05%  "
$ws.Range("D12").Value = "1.928.51"
$ws.Range("E12").Value = "  +1.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.044"
$ws.Range("E13").Value = "  +3.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.209"
$ws.Range("E14").Value = "  +2.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.09"
$ws.Range("E15").Value = "  +2.05%  "
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.007"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.06761"
$ws.Range("E17").Value = "  +2.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001039"
$ws.Range("E18").Value = "  +1.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.75"
$ws.Range("E19").Value = "  +1.35%  "
$ws.Range("E20").Value = "  +0.43%  "
$ws.Range("D21").Value = "29.510.72"
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.632"
$ws.Range("E22").Value = "  +2.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.74"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.177"
$ws.Range("E24").Value = "  -1.09%  "
$ws.Range("D25").Value = "2.150.02"
$ws.Range("E25").Value = "  +1.47%  "
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.567"
$ws.Range("E26").Value = "  +9.05%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.45"
$ws.Range("E27").Value = "  +1.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.10"
$ws.Range("E28").Value = "  +2.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.125"
$ws.Range("E29").Value = "  +1.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.62"
$ws.Range("E30").Value = "  +2.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.028"
$ws.Range("E31").Value = "  -3.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09545"
$ws.Range("E32").Value = "  +0.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.524"
$ws.Range("E33").Value = "  +3.20%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.560"
$ws.Range("E34").Value = "  +0.21%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.393"
$ws.Range("E35").Value = "  -1.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02279"
$ws.Range("E36").Value = "  +1.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06126"
$ws.Range("E37").Value = "  +0.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.187"
$ws.Range("E38").Value = "  +1.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.92"
$ws.Range("E39").Value = "  +8.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5977"
$ws.Range("E40").Value = "  +2.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.049"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1856"
$ws.Range("E42").Value = "  +1.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.408"
$ws.Range("E43").Value = "  -3.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.278"
$ws.Range("E44").Value = "  +0.51%  "
$ws.Range("E45").Value = "  +3.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07619"
$ws.Range("E46").Value = "  -0.97%  "
$ws.Range("E47").Value = "  +1.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.955"
$ws.Range("E48").Value = "  +2.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "116.61"
$ws.Range("E49").Value = "  +2.99%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.77"
$ws.Range("E50").Value = "  +2.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.417"
$ws.Range("E51").Value = "  +3.22%  "
